$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. New row 20: Pex tube 20x2 isole rouge (enter B20/E20 text first so the shared
#    string table grows in the same order a human filling the new row, then fixing
#    the "taille" column, then adding the "disponible" column would produce).
$ws.Range("A20").Value = "tube-alpex/Tube-Alpex-isolé-rouge.png"

$ws.Range("B20").Value = "Pex tube 20x2 isolé rouge"
$ws.Range("C2").Copy() | Out-Null
$ws.Range("B20").PasteSpecial(-4122) | Out-Null

$ws.Range("C20").Value = "20-50m"

$ws.Range("D20").Value = 158.75
$ws.Range("D19").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4122) | Out-Null

$ws.Range("E20").Value = "43123020"

# 2. Replace the generic "16"/"20" taille labels with the more specific
#    "16-25m" / "16-50m" / "20-25m" / "20-50m" values.
$ws.Range("C2").Value = "16-25m"
$ws.Range("C3").Value = "16-25m"
$ws.Range("C4").Value = "16-50m"
$ws.Range("C5").Value = "16-50m"
$ws.Range("C6").Value = "20-25m"
$ws.Range("C7").Value = "20-25m"
$ws.Range("C13").Value = "20-50m"
$ws.Range("C14").Value = "20-50m"
$ws.Range("C15").Value = "16-25m"
$ws.Range("C16").Value = "16-50m"
$ws.Range("C17").Value = "16-25m"
$ws.Range("C18").Value = "16-50m"
$ws.Range("C19").Value = "20-50m"

# 3. New "disponible" column (F) marking every article as available.
$ws.Range("F1").Value = "disponible"
$ws.Range("A1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null

for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 6).Value = "1"
}
$ws.Range("E2").Copy() | Out-Null
$ws.Range("F2:F20").PasteSpecial(-4122) | Out-Null

# 4. Selection left where the editor last clicked after finishing the update.
$ws.Range("A2").Select() | Out-Null
